$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data was recorded for "Vega Modelo de Temuco" (Acelga).
# It belongs chronologically between the existing row 168 (date 44217)
# and row 169 (date 44326), so insert a fresh row 169 and push every
# row below it (169-268) down by one -- the former row 268 becomes 269.
$ws.Rows.Item(169).Insert()

$ws.Range("A169").Value = 10
$ws.Range("B169").Value = "Vega Modelo de Temuco"
$ws.Range("C169").Value = "La Araucanía"
$ws.Range("D169").Value = 44596
$ws.Range("E169").Value = 9
$ws.Range("F169").Value = 100112009
$ws.Range("G169").Value = "Acelga"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 70
$ws.Range("K169").Value = 7000
$ws.Range("L169").Value = 8000
$ws.Range("M169").Value = 7429
$ws.Range("N169").Value = "$/docena de atados (12 kilos)"
$ws.Range("O169").Value = "Provincia de Cautín"
$ws.Range("P169").Value = 619
$ws.Range("Q169").Value = 12
$ws.Range("R169").Value = "Hortaliza"
